$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.369.14"
$ws.Range("E2").Value = "  +0.78%  "

# Row 3
$ws.Range("D3").Value = "2.236.08"
$ws.Range("E3").Value = "  -0.29%  "

# Row 4
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.38"
$ws.Range("E5").Value = "  -0.63%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.628"
$ws.Range("E6").Value = "  -0.78%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.52"
$ws.Range("E7").Value = "  -3.62%  "

# Row 8
$ws.Range("E8").Value = "  -0.08%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.617"
$ws.Range("E9").Value = "  -1.40%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.31"
$ws.Range("E10").Value = "  +4.03%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0948"
$ws.Range("E11").Value = "  -0.85%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.13"
$ws.Range("E12").Value = "  -0.47%  "

# Row 13
$ws.Range("E13").Value = "  -0.18%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.44"
$ws.Range("E14").Value = "  -2.54%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.852"
$ws.Range("E15").Value = "  -1.20%  "

# Row 16
$ws.Range("D16").Value = "2.239.61"
$ws.Range("E16").Value = "  -0.26%  "

# Row 17
$ws.Range("D17").Value = "42.198.06"
$ws.Range("E17").Value = "  +0.48%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000105"
$ws.Range("E18").Value = "  +6.74%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.15"
$ws.Range("E19").Value = "  +0.43%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.92"
$ws.Range("E20").Value = "  +0.01%  "

# Row 21
$ws.Range("E21").Value = "  +39.74%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.65"
$ws.Range("E22").Value = "  -0.49%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.18"
$ws.Range("E23").Value = "  -6.31%  "

# Row 24
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.58"
$ws.Range("E24").Value = "  +1.73%  "

# Row 25
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.01%  "

# Row 26
$ws.Range("E26").Value = "  +0.28%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.30"
$ws.Range("E27").Value = "  +0.22%  "

# Row 28
$ws.Range("E28").Value = "  +3.70%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.51"
$ws.Range("E29").Value = "  -2.10%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.92"
$ws.Range("E30").Value = "  +1.83%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.87"
$ws.Range("E31").Value = "  +19.80%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0808"
$ws.Range("E32").Value = "  -2.42%  "

# Row 33
$ws.Range("E33").Value = "  -3.23%  "

# Row 34
$ws.Range("E34").Value = "  -0.34%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "29.30"
$ws.Range("E35").Value = "  -13.24%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.54"
$ws.Range("E36").Value = "  +0.40%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0309"
$ws.Range("E37").Value = "  +2.38%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "13.19"
$ws.Range("E38").Value = "  -8.29%  "

# Row 39
$ws.Range("E39").Value = "  -1.00%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.65"
$ws.Range("E40").Value = "  -4.77%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "63.43"
$ws.Range("E41").Value = "  +4.26%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.200"
$ws.Range("E42").Value = "  -1.47%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.81"
$ws.Range("E43").Value = "  +1.44%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "105.77"
$ws.Range("E44").Value = "  -6.30%  "

# Row 45
$ws.Range("E45").Value = "  +2.11%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.994"
$ws.Range("E46").Value = "  -0.47%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.39"
$ws.Range("E47").Value = "  +4.45%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.14"
$ws.Range("E48").Value = "  +0.46%  "

# Row 49
$ws.Range("E49").Value = "  +0.17%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.73"
$ws.Range("E50").Value = "  +1.37%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.10"
$ws.Range("E51").Value = "  -2.51%  "
